$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Joe" to "Joseph" throughout column B (rows 3, 7, 11, 15, 19)
$ws.Range("B3").Value = 'Joseph'
$ws.Range("B7").Value = 'Joseph'
$ws.Range("B11").Value = 'Joseph'
$ws.Range("B15").Value = 'Joseph'
$ws.Range("B19").Value = 'Joseph'

# Fill in new sprint log entries for meeting dates 3/4/2019 (col F) and 3/6/2019 (col G)
$ws.Range("F2").Value = 'Not enough time to complete anything. Had other work.'
$ws.Range("G2").Value = 'Completed basic design for forgot information page. Completed back end for forgot information page however does not consider user account or databse, will have to be updated for security question answer verification'
$ws.Range("F3").Value = 'I created maze concept for level 4'
$ws.Range("G3").Value = 'I worked on the security issue where the user I signed out after 5 minuts of activity'
$ws.Range("F4").Value = 'No significant progress'
$ws.Range("G4").Value = 'No significant progress'
$ws.Range("F5").Value = 'Continued sprit work and started layout for login page'
$ws.Range("G5").Value = 'Continued sprite work and completed the layout for the login screen'
$ws.Range("F6").Value = 'Creating more forms in the user settings'
$ws.Range("G6").Value = 'Will spend time looking into firebase utiliaztion'
$ws.Range("F7").Value = 'I will work on my assigned issues'
$ws.Range("G7").Value = 'I will work on my assigned issues. Specifically the same issue as before. '
$ws.Range("F8").Value = 'Work on my assignments'
$ws.Range("G8").Value = 'Further work on assignments'
$ws.Range("F9").Value = 'Continue sprite work and complete login page'
$ws.Range("G9").Value = 'Continue sprite work and start coding the login screen'
$ws.Range("F10").Value = 'Not currently'
$ws.Range("G10").Value = 'Not currently'
$ws.Range("F11").Value = 'Errand for spring break'
$ws.Range("G11").Value = 'Errands for spring break'
$ws.Range("F12").Value = 'London'
$ws.Range("G12").Value = 'London'
$ws.Range("F13").Value = 'Nothing is currently getting in the way of my work'
$ws.Range("G13").Value = 'Nothing is currently getting in the way of my work'
$ws.Range("F14").Value = 'Nothing since last meeting'
$ws.Range("G14").Value = 'Nothing new as of yet'
$ws.Range("F15").Value = 'Documentation is tricky to master'
$ws.Range("G15").Value = 'Consistency is hard'
$ws.Range("F16").Value = 'Nothing new'
$ws.Range("G16").Value = 'Nothing new'
$ws.Range("F17").Value = 'Learned a little about creating more complex page layouts'
$ws.Range("G17").Value = 'Began to learn how to create a login screen'
$ws.Range("F18").Value = 'Not currently'
$ws.Range("G18").Value = 'Not currently'
$ws.Range("F19").Value = 'No'
$ws.Range("G19").Value = 'No'
$ws.Range("F20").Value = 'Not yet'
$ws.Range("G20").Value = 'Not yet'
$ws.Range("F21").Value = 'No changes need to be made to the project currently'
$ws.Range("G21").Value = 'No changes need to be made to the project currently'

# Update the saved view state (scroll position / active selection)
$ws.Range("D19").Select()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("F20").Select()
